$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new data rows (39 and 40) after the existing data (ending at row 38)
$ws.Range("A39").Value = 0.596596
$ws.Range("B39").Value = 0.644644
$ws.Range("C39").Value = 0.05961267486908179
$ws.Range("D39").Value = "query"

$ws.Range("A40").Value = 0.596596
$ws.Range("B40").Value = 0.627627
$ws.Range("C40").Value = 0.1126700295757865
$ws.Range("D40").Value = "query"
